$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 78 - this shifts the existing rows 78.. down by
# one and adjusts relative formula references the way Excel does.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row with the "Ford Transit Mk3 Minibus" entry.
$ws.Cells.Item(78, 1).Value = "Ford Transit Mk3 Minibus"
$ws.Cells.Item(78, 2).Value = 1986
$ws.Cells.Item(78, 3).Formula = "=B78-B77"
$ws.Cells.Item(78, 4).Value = 1
$ws.Cells.Item(78, 5).Value = "Car"
$ws.Cells.Item(78, 6).Formula = "=IF(B78 > 1900, ((B78-1900)*10)+400+D78, ((B78-1730)*2)+D78)+VLOOKUP(E78,'ID Scheme'!`$A`$2:`$B`$6,2, FALSE)"
$ws.Cells.Item(78, 7).Value = 80
$ws.Cells.Item(78, 8).Value = 9
$ws.Cells.Item(78, 9).Formula = "=SQRT(G78*H78)/`$B`$1"
$ws.Cells.Item(78, 10).Formula = "=I78*0.9"
$ws.Cells.Item(78, 11).Value = "x"
$ws.Cells.Item(78, 15).Formula = "=CONCATENATE(ROUND(L78*VLOOKUP(E78,'ID Scheme'!`$A`$2:`$E`$5,3),0), ""x"",ROUND(M78*VLOOKUP(E78,'ID Scheme'!`$A`$2:`$E`$5,5),0), ""x"",ROUND(N78*VLOOKUP(E78,'ID Scheme'!`$A`$2:`$E`$5,4),0))"

# The row-insert leaves the "Gap to Previous" formulas in the rows below the
# insertion point pointing at their pre-insert neighbour (row 77) instead of
# the row immediately above them post-insert. Re-fill them so they again read
# "this row's Intro Year minus the row above".
for ($r = 79; $r -le 82; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 3).Formula = "=B$r-B$prev"
}

$wb.Save()
